$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (2021) below the existing 2011-2020 table,
# carrying over the same formatting used by the preceding year rows.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 258411
